$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Bitcoin"
$ws.Range("C2").Value = "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.397.99"
$ws.Range("E2").Value = "  +1.07%  "
$ws.Range("B3").Value = "Ethereum"
$ws.Range("C3").Value = "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.010.76"
$ws.Range("E3").Value = "  +4.85%  "
$ws.Range("B4").Value = "TetherUSD"
$ws.Range("C4").Value = "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.96"
$ws.Range("E5").Value = "  +1.36%  "
$ws.Range("B6").Value = "USDC"
$ws.Range("C6").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("B7").Value = "XRP"
$ws.Range("C7").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5106"
$ws.Range("E7").Value = "  +1.20%  "
$ws.Range("B8").Value = "Cardano"
$ws.Range("C8").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4158"
$ws.Range("E8").Value = "  +3.31%  "
$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08785"
$ws.Range("E9").Value = "  +6.45%  "
$ws.Range("B10").Value = "Polygon"
$ws.Range("C10").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.136"
$ws.Range("E10").Value = "  +2.47%  "
$ws.Range("B11").Value = "OKB"
$ws.Range("C11").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.92"
$ws.Range("E11").Value = "  +2.10%  "
$ws.Range("B12").Value = "Solana"
$ws.Range("C12").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.67"
$ws.Range("E12").Value = "  +2.53%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.003.34"
$ws.Range("E13").Value = "  +4.81%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.605"
$ws.Range("E14").Value = "  +2.71%  "
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.485"
$ws.Range("E15").Value = "  +2.64%  "
$ws.Range("B16").Value = "BinanceUSD"
$ws.Range("C16").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.004"
$ws.Range("E16").Value = "  +0.41%  "
$ws.Range("B17").Value = "Litecoin"
$ws.Range("C17").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "94.32"
$ws.Range("E17").Value = "  +2.25%  "
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001117"
$ws.Range("E18").Value = "  +1.78%  "
$ws.Range("B19").Value = "TRON"
$ws.Range("C19").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06531"
$ws.Range("E19").Value = "  +0.35%  "
$ws.Range("B20").Value = "Avalanche"
$ws.Range("C20").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.97"
$ws.Range("E20").Value = "  +4.81%  "
$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.000"
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.230"
$ws.Range("E22").Value = "  +4.68%  "
$ws.Range("B23").Value = "WrappedBTC"
$ws.Range("C23").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.464.80"
$ws.Range("E23").Value = "  +1.17%  "
$ws.Range("B24").Value = "Cosmos"
$ws.Range("C24").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.98"
$ws.Range("E24").Value = "  +5.92%  "
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.226"
$ws.Range("E25").Value = "  +1.21%  "
$ws.Range("B26").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C26").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.236.51"
$ws.Range("E26").Value = "  +4.82%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.36"
$ws.Range("E27").Value = "  -0.47%  "
$ws.Range("B28").Value = "Monero"
$ws.Range("C28").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "162.98"
$ws.Range("E28").Value = "  +0.73%  "
$ws.Range("B29").Value = "LidoDAOToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.422"
$ws.Range("E29").Value = "  +6.19%  "
$ws.Range("B30").Value = "BitcoinCash"
$ws.Range("C30").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "131.51"
$ws.Range("E30").Value = "  +1.84%  "
$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.138"
$ws.Range("E31").Value = "  +0.89%  "
$ws.Range("B32").Value = "Stellar"
$ws.Range("C32").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1052"
$ws.Range("E32").Value = "  +1.35%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.119"
$ws.Range("E33").Value = "  +1.86%  "
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.831"
$ws.Range("E34").Value = "  +1.22%  "
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.351"
$ws.Range("E35").Value = "  +12.50%  "
$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02522"
$ws.Range("E36").Value = "  +3.07%  "
$ws.Range("B37").Value = "InternetComputer(DFINITY)"
$ws.Range("C37").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.466"
$ws.Range("E37").Value = "  +2.07%  "
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06615"
$ws.Range("E38").Value = "  +2.93%  "
$ws.Range("B39").Value = "Aptos"
$ws.Range("C39").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "12.40"
$ws.Range("E39").Value = "  +8.41%  "
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "9.122"
$ws.Range("E40").Value = "  +3.98%  "
$ws.Range("B41").Value = "Algorand"
$ws.Range("C41").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2197"
$ws.Range("E41").Value = "  +1.55%  "
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6670"
$ws.Range("E42").Value = "  +2.11%  "
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.233"
$ws.Range("E43").Value = "  +0.86%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.63"
$ws.Range("E44").Value = "  +2.53%  "
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6177"
$ws.Range("E45").Value = "  +2.30%  "
$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.199"
$ws.Range("E46").Value = "  -0.04%  "
$ws.Range("B47").Value = "PancakeSwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.672"
$ws.Range("E47").Value = "  +0.90%  "
$ws.Range("B48").Value = "EOS"
$ws.Range("C48").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.272"
$ws.Range("E48").Value = "  +4.71%  "
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "124.32"
$ws.Range("E49").Value = "  +0.52%  "
$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "81.30"
$ws.Range("E50").Value = "  +3.32%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06896"
$ws.Range("E51").Value = "  +1.51%  "
